$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1350
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1350
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1350
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -1700
$ws.Range("H116").Value = 593910.9
$ws.Range("J116").Value = 8419.6
$ws.Range("L116").Value = 8419.6
$ws.Range("N116").Value = -15303.6
$ws.Range("H123").Value = 42597.145
$ws.Range("J123").Value = 42597.145
$ws.Range("L123").Value = 42597.145
$ws.Range("N123").Value = -52397.145
$ws.Range("H129").Value = 959.8591300000001
$ws.Range("J129").Value = 990.34326
$ws.Range("L129").Value = 2971.02978
$ws.Range("N129").Value = -12971.02978
$ws.Range("H137").Value = 2341.9272
$ws.Range("I137").Value = 1451.9111
$ws.Range("J137").Value = 6347
$ws.Range("K137").Value = 4355.7333
$ws.Range("L137").Value = 19041
$ws.Range("M137").Value = -1805.7333
$ws.Range("N137").Value = -24141
$ws.Range("H138").Value = 3137.0928
$ws.Range("I138").Value = 1632.7778
$ws.Range("J138").Value = 3479.8481
$ws.Range("K138").Value = 4898.3334
$ws.Range("L138").Value = 10439.5443
$ws.Range("M138").Value = 241.6665999999996
$ws.Range("N138").Value = -20719.5443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6298.9
$ws.Range("I32").Value = 3646.1333
$ws.Range("J32").Value = 14257.2
$ws.Range("K32").Value = 3646.1333
$ws.Range("L32").Value = 14257.2
$ws.Range("M32").Value = -3359.1333
$ws.Range("N32").Value = -14831.2
$ws.Range("H45").Value = 1623.0526
$ws.Range("I45").Value = 819.8333
$ws.Range("K45").Value = 819.8333
$ws.Range("M45").Value = -442.8333
$ws.Range("H61").Value = 1945.1538
$ws.Range("I61").Value = 1455.238
$ws.Range("J61").Value = 4002.8
$ws.Range("K61").Value = 1455.238
$ws.Range("L61").Value = 4002.8
$ws.Range("M61").Value = -1243.238
$ws.Range("N61").Value = -4426.8
$ws.Range("H74").Value = 1316.8182
$ws.Range("I74").Value = 846.82355
$ws.Range("J74").Value = 2914.8
$ws.Range("K74").Value = 846.82355
$ws.Range("L74").Value = 2914.8
$ws.Range("M74").Value = 27.17645000000005
$ws.Range("N74").Value = -4662.8
$ws.Range("H77").Value = 1316.8182
$ws.Range("I77").Value = 846.82355
$ws.Range("J77").Value = 2914.8
$ws.Range("K77").Value = 4234.117749999999
$ws.Range("L77").Value = 14574
$ws.Range("M77").Value = 133.8822500000006
$ws.Range("N77").Value = -23310
$ws.Range("H110").Value = 397.33334
$ws.Range("I110").Value = 367.2
$ws.Range("K110").Value = 367.2
$ws.Range("M110").Value = 1677.8
$ws.Range("H132").Value = 2667.1714
$ws.Range("I132").Value = 1348.8182
$ws.Range("K132").Value = 4046.4546
$ws.Range("M132").Value = -1516.4546
$ws.Range("H136").Value = 1945.1538
$ws.Range("I136").Value = 1455.238
$ws.Range("J136").Value = 4002.8
$ws.Range("K136").Value = 4365.714
$ws.Range("L136").Value = 12008.4
$ws.Range("M136").Value = -1815.714
$ws.Range("N136").Value = -17108.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = $null
$ws.Range("H134").Value = 2915.491
$ws.Range("I134").Value = 1548.525
$ws.Range("J134").Value = 6560.7334
$ws.Range("K134").Value = 4645.575000000001
$ws.Range("L134").Value = 19682.2002
$ws.Range("M134").Value = -2110.575000000001
$ws.Range("N134").Value = -24752.2002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4136.25
$ws.Range("I31").Value = 1781.1765
$ws.Range("J31").Value = 6805.3335
$ws.Range("K31").Value = 1781.1765
$ws.Range("L31").Value = 6805.3335
$ws.Range("M31").Value = -1486.1765
$ws.Range("N31").Value = -7395.3335
$ws.Range("H34").Value = 4136.25
$ws.Range("I34").Value = 1781.1765
$ws.Range("J34").Value = 6805.3335
$ws.Range("K34").Value = 1781.1765
$ws.Range("L34").Value = 6805.3335
$ws.Range("M34").Value = -1579.1765
$ws.Range("N34").Value = -7209.3335
$ws.Range("H58").Value = 1907.0476
$ws.Range("I58").Value = 1658.0984
$ws.Range("J58").Value = 9500
$ws.Range("K58").Value = 1658.0984
$ws.Range("L58").Value = 9500
$ws.Range("M58").Value = -1455.0984
$ws.Range("N58").Value = -9906
$ws.Range("H109").Value = 29998.545
$ws.Range("J109").Value = 29998.545
$ws.Range("L109").Value = 29998.545
$ws.Range("N109").Value = -32078.545
$ws.Range("H132").Value = 3306.647
$ws.Range("I132").Value = 2871.077
$ws.Range("J132").Value = 4722.25
$ws.Range("K132").Value = 8613.231
$ws.Range("L132").Value = 14166.75
$ws.Range("M132").Value = -6083.231
$ws.Range("N132").Value = -19226.75
$ws.Range("H134").Value = 3328.25
$ws.Range("I134").Value = 3470.3555
$ws.Range("J134").Value = 2901.9333
$ws.Range("K134").Value = 10411.0665
$ws.Range("L134").Value = 8705.7999
$ws.Range("M134").Value = -7876.066500000001
$ws.Range("N134").Value = -13775.7999
$ws.Range("H136").Value = 1907.0476
$ws.Range("I136").Value = 1658.0984
$ws.Range("J136").Value = 9500
$ws.Range("K136").Value = 4974.2952
$ws.Range("L136").Value = 28500
$ws.Range("M136").Value = -2424.2952
$ws.Range("N136").Value = -33600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2969.8125
$ws.Range("I5").Value = 875
$ws.Range("J5").Value = 3269.0715
$ws.Range("K5").Value = 2625
$ws.Range("L5").Value = 9807.2145
$ws.Range("M5").Value = -2513
$ws.Range("N5").Value = -10031.2145
$ws.Range("H106").Value = 4000
$ws.Range("J106").Value = 4000
$ws.Range("L106").Value = 12000
$ws.Range("N106").Value = -13892
$ws.Range("H133").Value = 3824.2856
$ws.Range("I133").Value = 3695.8333
$ws.Range("J133").Value = 3995.5557
$ws.Range("K133").Value = 11087.4999
$ws.Range("L133").Value = 11986.6671
$ws.Range("M133").Value = -6027.499899999999
$ws.Range("N133").Value = -22106.6671
$ws.Range("H134").Value = 3456.2222
$ws.Range("I134").Value = 2314.5334
$ws.Range("J134").Value = 4883.3335
$ws.Range("K134").Value = 6943.600199999999
$ws.Range("L134").Value = 14650.0005
$ws.Range("M134").Value = -1873.600199999999
$ws.Range("N134").Value = -24790.0005
$ws.Range("H135").Value = 2969.8125
$ws.Range("I135").Value = 875
$ws.Range("J135").Value = 3269.0715
$ws.Range("K135").Value = 7875
$ws.Range("L135").Value = 29421.6435
$ws.Range("M135").Value = -5340
$ws.Range("N135").Value = -34491.6435
$ws.Range("H136").Value = 3293.6365
$ws.Range("I136").Value = 2947.1428
$ws.Range("J136").Value = 3900
$ws.Range("K136").Value = 8841.428400000001
$ws.Range("L136").Value = 11700
$ws.Range("M136").Value = -3741.428400000001
$ws.Range("N136").Value = -21900
$ws.Range("H137").Value = 3075.625
$ws.Range("I137").Value = 2829.2856
$ws.Range("J137").Value = 4800
$ws.Range("K137").Value = 8487.856800000001
$ws.Range("L137").Value = 14400
$ws.Range("M137").Value = -3387.856800000001
$ws.Range("N137").Value = -24600
$ws.Range("H138").Value = 2880.2083
$ws.Range("I138").Value = 1564.5385
$ws.Range("J138").Value = 4435.091
$ws.Range("K138").Value = 4693.6155
$ws.Range("L138").Value = 13305.273
$ws.Range("M138").Value = 446.3845000000001
$ws.Range("N138").Value = -23585.273
$ws.Range("H139").Value = 3681.5908
$ws.Range("I139").Value = 1999.75
$ws.Range("J139").Value = 8166.5
$ws.Range("K139").Value = 5999.25
$ws.Range("L139").Value = 24499.5
$ws.Range("M139").Value = -859.25
$ws.Range("N139").Value = -34779.5
$ws.Range("H140").Value = 1330.6154
$ws.Range("I140").Value = 1191.5
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 3574.5
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = 1605.5
$ws.Range("N140").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5235086
$ws.Range("I11").Value = 7500833.5
$ws.Range("K11").Value = 7500833.5
$ws.Range("M11").Value = -7500694.5
$ws.Range("H113").Value = 1337
$ws.Range("I113").Value = 1005.5
$ws.Range("K113").Value = 1005.5
$ws.Range("M113").Value = 1164.5
$ws.Range("H132").Value = 2323.4167
$ws.Range("I132").Value = 926.06665
$ws.Range("J132").Value = 2958.5757
$ws.Range("K132").Value = 2778.19995
$ws.Range("L132").Value = 8875.7271
$ws.Range("M132").Value = -248.1999500000002
$ws.Range("N132").Value = -13935.7271

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 625.8946999999999
$ws.Range("I16").Value = 625.8946999999999
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 625.8946999999999
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -455.8946999999999
$ws.Range("N16").Value = $null
$ws.Range("H61").Value = 1742.091
$ws.Range("I61").Value = 1616.3
$ws.Range("K61").Value = 1616.3
$ws.Range("M61").Value = -1414.3
$ws.Range("H108").Value = 39800
$ws.Range("J108").Value = 39800
$ws.Range("L108").Value = 39800
$ws.Range("N108").Value = -47480
$ws.Range("H113").Value = 1742.091
$ws.Range("I113").Value = 1616.3
$ws.Range("K113").Value = 1616.3
$ws.Range("M113").Value = 553.7
$ws.Range("H114").Value = 41666.668
$ws.Range("J114").Value = 41666.668
$ws.Range("L114").Value = 41666.668
$ws.Range("N114").Value = -50344.668
$ws.Range("H116").Value = 41400
$ws.Range("J116").Value = 41400
$ws.Range("L116").Value = 41400
$ws.Range("N116").Value = -50578
$ws.Range("H123").Value = 31042.9
$ws.Range("J123").Value = 31042.9
$ws.Range("L123").Value = 31042.9
$ws.Range("N123").Value = -40842.9
$ws.Range("H132").Value = 4782.385
$ws.Range("I132").Value = 1922.6154
$ws.Range("J132").Value = 7642.154
$ws.Range("K132").Value = 5767.8462
$ws.Range("L132").Value = 22926.462
$ws.Range("M132").Value = -3237.8462
$ws.Range("N132").Value = -27986.462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 27377
$ws.Range("J109").Value = 27377
$ws.Range("L109").Value = 27377
$ws.Range("N109").Value = -30151
$ws.Range("H132").Value = 14499831
$ws.Range("I132").Value = 15285.714
$ws.Range("K132").Value = 45857.142
$ws.Range("M132").Value = -43327.142
